# Apply the edit described by the diff:
#  - On sheet "Planilha2", cell B2 text changes from "Kingston Mouse" to "KingstonMouse"
#  - The active selection on "Planilha2" changes from A2 to E7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha2")

# Update cell B2 value (removing the space from "Kingston Mouse")
$ws.Range("B2").Value = "KingstonMouse"

# Activate the sheet and move/update the selection to E7
$ws.Activate()
$ws.Range("E7").Select() | Out-Null
